# tc2_buyer1.xlsx - "calc fix first tries"
# Insert a new data row (supplier1 / 11dd / 10000 / 28-Sep-2016) between the
# existing two rows, fix up the surrounding values, and touch up the
# formatting that Excel re-saved along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 2 - shifts old row 2 down to row 3.
$ws.Rows.Item(2).Insert()

# Row 1: D1 was stored as the shared string "28.08.2016"; it becomes a real
# date value (serial 42641 == 28-Sep-2016) formatted with the existing date
# style already applied to the column.
$ws.Range("D1").Value = 42641

# Row 2 (brand new row): supplier1 / 11dd / 10000 / same date serial as D1.
$ws.Range("A2").Value = "11dd"
$ws.Range("B2").Value = "supplier1"
$ws.Range("C2").Value = 10000
$ws.Range("D2").Value = 42641

# Row 3 (previously row 2): the amount changes from 3000 to 20000.
$ws.Range("C3").Value = 20000

# Cosmetic follow-up matching the resave: widen column D and move the
# active selection (both happened as part of the same Excel session).
$ws.Columns.Item(4).ColumnWidth = 9.25
[void]$ws.Range("D11").Select()

Write-Output "edits applied"
